$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("B2").Value = 0.00231176578588815
$ws.Range("C2").Value = 0.488201967721922
$ws.Range("B3").Value = 0.00235600679767414
$ws.Range("C3").Value = 0.60088414031413
$ws.Range("B4").Value = 0.00316367984963938
$ws.Range("C4").Value = 0.408652337854738
$ws.Range("B5").Value = 0.00257243083110167
$ws.Range("C5").Value = 0.641599854107281
$ws.Range("B6").Value = 0.00308566566691175
$ws.Range("C6").Value = 0.584669548600328
$ws.Range("B7").Value = 0.00261840705503407
$ws.Range("C7").Value = 0.489618402585804
$ws.Range("B8").Value = 0.00314253548686994
$ws.Range("C8").Value = 0.722996518942883
$ws.Range("B9").Value = 0.0023030426190356
$ws.Range("C9").Value = 0.718035372338152
$ws.Range("B10").Value = 0.00244571428180775
$ws.Range("C10").Value = 0.705526558186348
$ws.Range("B11").Value = 0.00232868099212124
$ws.Range("C11").Value = 0.719430138097316
$ws.Range("B12").Value = 0.00631395726223987
$ws.Range("C12").Value = 0.530861890919076
$ws.Range("B13").Value = 0.00319887393836746
$ws.Range("C13").Value = 0.720300709917248
$ws.Range("B14").Value = 0.00252262354389423
$ws.Range("C14").Value = 0.665732600684077
$ws.Range("B15").Value = 0.0020721419568032
$ws.Range("C15").Value = 0.790418598605476
$ws.Range("B16").Value = 0.00257101062139269
$ws.Range("C16").Value = 0.685866938415752
$ws.Range("B17").Value = 0.00322688386575513
$ws.Range("C17").Value = 0.648210177577909
$ws.Range("B18").Value = 0.00299701596993502
$ws.Range("C18").Value = 0.707861024672471
$ws.Range("B19").Value = 0.002805812016169
$ws.Range("C19").Value = 0.680251933957271
$ws.Range("B20").Value = 0.00331646182502219
$ws.Range("C20").Value = 0.556423239324126
$ws.Range("B21").Value = 0.0030954398158163
$ws.Range("C21").Value = 0.547911756027029
$ws.Range("B22").Value = 0.00326586006576449
$ws.Range("C22").Value = 0.58768053430424
$ws.Range("B23").Value = 0.00259054524987445
$ws.Range("C23").Value = 0.745087399601091
$ws.Range("B24").Value = 0.00303694450629423
$ws.Range("C24").Value = 0.626722143122576
$ws.Range("B25").Value = 0.0031033199016924
$ws.Range("C25").Value = 0.518445365220124
$ws.Range("B26").Value = 0.0348687935198066
$ws.Range("C26").Value = 0.741373240461102
$ws.Range("B27").Value = 0.0395065860548887
$ws.Range("C27").Value = 0.675041949588754
$ws.Range("B28").Value = 0.0346179289538455
$ws.Range("C28").Value = 0.740311213456729
$ws.Range("B29").Value = 0.0580035365813083
$ws.Range("C29").Value = 0.627805405405405
$ws.Range("B30").Value = 0.0320161469350603
$ws.Range("C30").Value = 0.77031886500913
$ws.Range("B31").Value = 0.0270950772561376
$ws.Range("C31").Value = 0.804945715031899
$ws.Range("B32").Value = 0.0367773481078161
$ws.Range("C32").Value = 0.709114589327502
$ws.Range("B33").Value = 0.0240508973905873
$ws.Range("C33").Value = 0.82915862527927
$ws.Range("B34").Value = 0.0259547617964963
$ws.Range("C34").Value = 0.84860080738909
